$d = $word.ActiveDocument

# Replace the trailing space in " Amar B " with "KC" to form " Amar B KC"
$d.Content.Find.Execute(" Amar B ", $true, $false, $false, $false, $false,
                         $true, 1, $false, " Amar B KC", 2)

# Remove the separate "kc" run that followed the bookmark
$d.Content.Find.Execute("kc", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
